$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current header row (B1:E1) and data row (B2:E2) values before overwriting anything
$h1 = $ws.Range("B1").Value()
$h2 = $ws.Range("C1").Value()
$h3 = $ws.Range("D1").Value()
$h4 = $ws.Range("E1").Value()

$d1 = $ws.Range("B2").Value()
$d2 = $ws.Range("C2").Value()
$d3 = $ws.Range("D2").Value()
$d4 = $ws.Range("E2").Value()

# Copy the header cell formatting (bold font + border + centered alignment) so it can be
# applied to the new A1 header cell
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift headers one column to the left: A1:D1 = old B1:E1
$ws.Range("A1").Value = $h1
$ws.Range("B1").Value = $h2
$ws.Range("C1").Value = $h3
$ws.Range("D1").Value = $h4

# Shift data row one column to the left: A2:D2 = old B2:E2 (unstyled, matches original B2:D2)
$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = $d1
$ws.Range("B2").Value = $d2
$ws.Range("C2").Value = $d3
$ws.Range("D2").Value = $d4

# Remove the now-unused column E entirely (content + formatting)
$ws.Range("E1:E2").Clear()
